# Add team record (Wins / Losses / Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new columns AD, AE, AF, formatted like the other
# header cells (copy formatting from the last existing header cell AC1).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-47: every row gets the same team record (88 wins, 74 losses, 0 ties).
for ($r = 2; $r -le 47; $r++) {
    $ws.Cells.Item($r, 30).Value = 88   # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 74   # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF -> Ties
}
